$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '242.78'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.09'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.414'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05892'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.437'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.530'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8087'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9321'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07370'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03295'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03059'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09365'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.854'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001573'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04674'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005901'
$ws.Range("E18").Value = '17OneONE'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005854'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001262'
$ws.Range("E20").Value = '19BitKanKANBestin24h'
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004900'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00006801'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.565'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.145'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3233'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1330'
$ws.Range("E26").Value = '25ProBitTokenPROB'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03972'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006190'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1072'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003000'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008747'
$ws.Range("E44").Value = '43LocalTradersLCT'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005189'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6701'
